# Final parent measurements added for KAHI participant 1 data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New "Post Treatment" (column D) responses for each question row (2-15).
$values = @{
    2  = "A little worse"
    3  = "Somewhat worse"
    4  = "Somewhat worse"
    5  = "Somewhat worse"
    6  = "Somewhat worse"
    7  = "A lot worse"
    8  = "A lot worse"
    9  = "Somewhat worse"
    10 = "A little worse"
    11 = "A lot worse"
    12 = "A lot worse"
    13 = "A lot worse"
    14 = "Somewhat worse"
    15 = "A lot worse"
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 4).Value = $values[$row]
}

# Widen column D slightly to fit the newly entered data.
$ws.Columns("D").ColumnWidth = 16.1

# Move the active selection down to the next empty row, as left by the author.
$ws.Range("D16").Select()
